$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.944.64'
$ws.Range('E2').Value = '  +1.65%  '
$ws.Range('D3').Value = '2.369.41'
$ws.Range('E3').Value = '  +0.74%  '
$ws.Range('E4').Value = '  +0.29%  '
$ws.Range('E5').Value = '  +6.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '240.99'
$ws.Range('E6').Value = '  +3.08%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '76.05'
$ws.Range('E7').Value = '  +6.67%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.622'
$ws.Range('E9').Value = '  +28.25%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.103'
$ws.Range('E10').Value = '  +5.00%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '57.47'
$ws.Range('E11').Value = '  +0.94%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '32.91'
$ws.Range('E12').Value = '  +20.65%  '
$ws.Range('E13').Value = '  +19.45%  '
$ws.Range('E14').Value = '  +1.53%  '
$ws.Range('D15').Value = '2.723.14'
$ws.Range('E15').Value = '  +0.43%  '
$ws.Range('E16').Value = '  +5.42%  '
$ws.Range('E17').Value = '  +6.67%  '
$ws.Range('D18').Value = '2.357.84'
$ws.Range('E18').Value = '  +0.16%  '
$ws.Range('D19').Value = '43.949.01'
$ws.Range('E19').Value = '  +1.68%  '
$ws.Range('E20').Value = '  +1.95%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.66'
$ws.Range('E21').Value = '  +5.37%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '77.98'
$ws.Range('E22').Value = '  +4.82%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '257.07'
$ws.Range('E23').Value = '  +2.86%  '
$ws.Range('E25').Value = '  +3.06%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.11'
$ws.Range('E26').Value = '  +10.91%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.63'
$ws.Range('E27').Value = '  -3.98%  '
$ws.Range('E28').Value = '  +17.81%  '
$ws.Range('E29').Value = '  +1.90%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '23.21'
$ws.Range('E30').Value = '  +3.85%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '175.54'
$ws.Range('E31').Value = '  +1.71%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.129'
$ws.Range('E32').Value = '  -1.00%  '
$ws.Range('E33').Value = '  +6.11%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.32'
$ws.Range('E34').Value = '  +6.97%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0756'
$ws.Range('E35').Value = '  +9.81%  '
$ws.Range('E36').Value = '  +5.78%  '
$ws.Range('E37').Value = '  +3.91%  '
$ws.Range('E38').Value = '  +1.93%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.51'
$ws.Range('E39').Value = '  -0.41%  '
$ws.Range('E40').Value = '  +7.73%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '18.97'
$ws.Range('E41').Value = '  +0.97%  '
$ws.Range('B42').Value = 'BinanceUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  +0.03%  '
$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.200'
$ws.Range('E43').Value = '  +19.11%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.91'
$ws.Range('E44').Value = '  +0.23%  '
$ws.Range('B45').Value = 'ARBITRUM'
$ws.Range('C45').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.20'
$ws.Range('E45').Value = '  +4.04%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.101'
$ws.Range('E46').Value = '  +5.12%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.52'
$ws.Range('E47').Value = '  +14.03%  '
$ws.Range('B48').Value = 'TrustWalletToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.26'
$ws.Range('E48').Value = '  +4.36%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '102.23'
$ws.Range('E49').Value = '  +3.25%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.51'
$ws.Range('E50').Value = '  +0.80%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '54.76'
$ws.Range('E51').Value = '  +8.93%  '
